$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row rework -------------------------------------------------
# Before: B1="Inicio de Conexi¢n", C1="Fin de Conexio", D1="Usuario"
# After:  B1="Usuario",            C1="Inicio de Conexi¢n"  (D1 removed)
$ws.Range("D1").Clear() | Out-Null
$ws.Range("C1").Value2 = "Inicio de Conexi¢n"
$ws.Range("B1").Value2 = "Usuario"

# --- New data rows 2..26 (A=id, B=user, C=connection start datetime) ---
$ids   = @(1,2,4,9,11,13,19,20,46,47,51,52,53,55,57,59,62,63,64,66,70,78,80,83,90)
$users = @("invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew","invitede-pansew")
$times = @(43705.42083333333,43705.42083333333,43705.42083333333,43705.42152777778,43705.42152777778,43705.42152777778,43705.42152777778,43705.42152777778,43705.42361111111,43705.42361111111,43705.42361111111,43705.42361111111,43705.42361111111,43705.42361111111,43705.42430555556,43705.42430555556,43705.42430555556,43705.42430555556,43705.42430555556,43705.425,43705.42569444444,43705.42638888889,43705.42638888889,43705.42708333334,43705.42777777778)

$count    = $ids.Length
$firstRow = 2
$lastRow  = $firstRow + $count - 1
$aRange   = "A" + $firstRow + ":A" + $lastRow
$cRange   = "C" + $firstRow + ":C" + $lastRow

# Write the raw values first.
for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 1).Value2 = $ids[$i]
    $ws.Cells.Item($r, 2).Value2 = $users[$i]
    $ws.Cells.Item($r, 3).Value2 = $times[$i]
}

# Column A reuses the same bold/centered/bordered style as the header
# cells (B1/C1) - copy-format it down instead of rebuilding it by hand.
$ws.Range("B1").Copy() | Out-Null
$ws.Range($aRange).PasteSpecial(-4122) | Out-Null

# Column C gets a custom date-time display format. Apply the lowercase
# form first, then the uppercase form to the anchor cell only, so the
# stylesheet ends up with numFmtId 164 (yyyy-mm-dd h:mm:ss) registered
# ahead of 165 (YYYY-MM-DD HH:MM:SS) - 165 being the format that is
# actually used. Then copy that one resulting style down the column so
# no extra cellXfs entries get minted.
$ws.Range("C" + $firstRow).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("C" + $firstRow).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C" + $firstRow).Copy() | Out-Null
$ws.Range($cRange).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# PasteSpecial(formats) only moves formatting - restore the values.
for ($i = 0; $i -lt $count; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 3).Value2 = $times[$i]
}

Write-Output "done"
